# timesheet.xlsx: add the 28/09 kanbanflow entry as row 9 and move the
# selection on to the next empty row (D10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dia / Hora Inicio / Hora Fim / Atividade for 28/09/2016, 21:00 - 22:20
$ws.Range("A9").Value = 42641
$ws.Range("A9").NumberFormat = "d-mmm"

$ws.Range("B9").Value = 0.875
$ws.Range("B9").NumberFormat = "h:mm"

$ws.Range("C9").Value = 0.93055555555555547
$ws.Range("C9").NumberFormat = "h:mm"

$ws.Range("D9").Value = "Implementação das tarefas (kanbanflow) do dia 28/09"

# Active cell moves down to the next row, same column, as in the source file.
$ws.Range("D10").Select() | Out-Null
